$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column F (Cargo) content left into column E, then delete the now-empty F.
$ws.Range("F1:F23").Cut($ws.Range("E1:E23"))
$ws.Columns("F").Delete()

# Restore the selection that Excel records for this sheet's last-saved view.
$ws.Range("E1:E1048576").Select() | Out-Null
